$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8736847402618082
$ws.Range("C2").Value = 0.4372017465268527
$ws.Range("E2").Value = 0.0801198944586119
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002446411695913761
$ws.Range("I2").Value = 1.496099152613581
$ws.Range("M2").Value = 0.4110552294794445
$ws.Range("B3").Value = 0.7868024024286342
$ws.Range("C3").Value = 0.3862978329114526
$ws.Range("E3").Value = 0.07491932498378873
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002452612398224435
$ws.Range("I3").Value = 1.415280720512087
$ws.Range("M3").Value = 0.3723175530131755
$ws.Range("B4").Value = 0.7340310585277336
$ws.Range("C4").Value = 0.3552235821381373
$ws.Range("E4").Value = 0.07177625357847006
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002456610528059306
$ws.Range("I4").Value = 1.366150344832135
$ws.Range("M4").Value = 0.3487910188779821
$ws.Range("B5").Value = 0.7126684126811824
$ws.Range("C5").Value = 0.3426040627689986
$ws.Range("E5").Value = 0.07050775215833482
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002458287986670606
$ws.Range("I5").Value = 1.346248331103439
$ws.Range("M5").Value = 0.3392673446715406
$ws.Range("B6").Value = 0.709129674874589
$ws.Range("C6").Value = 0.340511169630787
$ws.Range("E6").Value = 0.07029785501515562
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002458569443539507
$ws.Range("I6").Value = 1.342950664421821
$ws.Range("M6").Value = 0.33768974439743
$ws.Range("B7").Value = 0.7337423821909113
$ws.Range("C7").Value = 0.355053217438865
$ws.Range("E7").Value = 0.0717590965756969
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002456632955468247
$ws.Range("I7").Value = 1.365881463793642
$ws.Range("M7").Value = 0.3486623237256197
$ws.Range("B8").Value = 0.8436069130013948
$ws.Range("C8").Value = 0.4196112843541755
$ws.Range("E8").Value = 0.07831617363488874
$ws.Range("F8").Value = 0.4248636149813905
$ws.Range("G8").Value = 0.00244851019818329
$ws.Range("I8").Value = 1.468127815093595
$ws.Range("M8").Value = 0.3976439313101849
$ws.Range("B9").Value = 1.063734972815325
$ws.Range("C9").Value = 0.5477432154464736
$ws.Range("E9").Value = 0.09158515289335867
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002434087193570554
$ws.Range("I9").Value = 1.672766980121295
$ws.Range("M9").Value = 0.4958190087955074
$ws.Range("B10").Value = 1.228510368539844
$ws.Range("C10").Value = 0.6429627307967962
$ws.Range("E10").Value = 0.101603379208143
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002424396147721848
$ws.Range("I10").Value = 1.825963540465693
$ws.Range("M10").Value = 0.5693487202740641
$ws.Range("B11").Value = 1.304176487997779
$ws.Range("C11").Value = 0.6865479544622985
$ws.Range("E11").Value = 0.1062237291304839
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002420181406213231
$ws.Range("I11").Value = 1.896347045886642
$ws.Range("M11").Value = 0.6031274223404779
$ws.Range("B12").Value = 1.332934551158644
$ws.Range("C12").Value = 0.7030938083909746
$ws.Range("E12").Value = 0.1079827304919831
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002418613050855915
$ws.Range("I12").Value = 1.923104814885789
$ws.Range("M12").Value = 0.6159677940907926
$ws.Range("B13").Value = 1.326736271492791
$ws.Range("C13").Value = 0.699528498118184
$ws.Range("E13").Value = 0.1076034760457887
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002418949596317262
$ws.Range("I13").Value = 1.917337292535677
$ws.Range("M13").Value = 0.6132001785232433
$ws.Range("B14").Value = 1.306540310389778
$ws.Range("C14").Value = 0.6879083550408609
$ws.Range("E14").Value = 0.1063682535170827
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002420051823152462
$ws.Range("I14").Value = 1.898546287346704
$ws.Range("M14").Value = 0.6041828132983795
$ws.Range("B15").Value = 1.294183458666623
$ws.Range("C15").Value = 0.6807961013576005
$ws.Range("E15").Value = 0.1056128744049332
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002420730567542939
$ws.Range("I15").Value = 1.887050095279818
$ws.Range("M15").Value = 0.5986658636297193
$ws.Range("B16").Value = 1.223579930016058
$ws.Range("C16").Value = 0.6401199534735724
$ws.Range("E16").Value = 0.1013027231274961
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002424675473796873
$ws.Range("I16").Value = 1.82137822389376
$ws.Range("M16").Value = 0.567147981205494
$ws.Range("B17").Value = 1.180450638835396
$ws.Range("C17").Value = 0.6152370995096135
$ws.Range("E17").Value = 0.0986749449421751
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002427145040436911
$ws.Range("I17").Value = 1.781272079656731
$ws.Range("M17").Value = 0.5478984171681276
$ws.Range("B18").Value = 1.155710357583416
$ws.Range("C18").Value = 0.6009503534518785
$ws.Range("E18").Value = 0.09716942793686201
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002428583718906017
$ws.Range("I18").Value = 1.758269051778086
$ws.Range("M18").Value = 0.5368575025631088
$ws.Range("B19").Value = 1.147345063009482
$ws.Range("C19").Value = 0.5961173667021171
$ws.Range("E19").Value = 0.09666069156364188
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002429073970874314
$ws.Range("I19").Value = 1.750491618447313
$ws.Range("M19").Value = 0.5331244987913522
$ws.Range("B20").Value = 1.185034916799225
$ws.Range("C20").Value = 0.6178832955160942
$ws.Range("E20").Value = 0.09895406233707149
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002426880263640379
$ws.Range("I20").Value = 1.785534685126265
$ws.Range("M20").Value = 0.5499443532297761
$ws.Range("B21").Value = 1.312469481172627
$ws.Range("C21").Value = 0.6913203390873832
$ws.Range("E21").Value = 0.106730811342338
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002419727322877718
$ws.Range("I21").Value = 1.904062767004859
$ws.Range("M21").Value = 0.6068300854230699
$ws.Range("B22").Value = 1.396368498237962
$ws.Range("C22").Value = 0.7395562323401919
$ws.Range("E22").Value = 0.1118681619062372
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002415213689251239
$ws.Range("I22").Value = 1.98214261826331
$ws.Range("M22").Value = 0.6442952728001927
$ws.Range("B23").Value = 1.351532894831507
$ws.Range("C23").Value = 0.7137890337996282
$ws.Range("E23").Value = 0.1091211419810207
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002417608010263004
$ws.Range("I23").Value = 1.940411896047152
$ws.Range("M23").Value = 0.6242725562362921
$ws.Range("B24").Value = 1.182962192638627
$ws.Range("C24").Value = 0.6166868921119431
$ws.Range("E24").Value = 0.09882785710345132
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002426999910341937
$ws.Range("I24").Value = 1.783607391750536
$ws.Range("M24").Value = 0.5490193049415808
$ws.Range("B25").Value = 1.003663527101594
$ws.Range("C25").Value = 0.5129010031465668
$ws.Range("E25").Value = 0.0879496205680681
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002437829082954392
$ws.Range("I25").Value = 1.616930341284188
$ws.Range("M25").Value = 0.469021949302757
